# Rekap Tamu — add "manage data divisi" API rows + refresh the report period
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash a clean copy of the existing data-row format (border + general
# number format, style index 2) far below the used range so we can keep
# reusing it as a formats-only paste source without it drifting as we
# overwrite cells' number formats below.
$ws.Range("A7:H7").Copy()
$ws.Range("A20:H20").PasteSpecial(-4122)

function Set-GuestRow($row, $no, $nama, $instansi, $telp, $tanggal, $jam, $jumlah, $bagian) {
    $ws.Cells.Item($row, 1).Value = $no
    $ws.Cells.Item($row, 2).Value = $nama
    $ws.Cells.Item($row, 3).Value = $instansi

    # Telepon / Tanggal / jam are kept as text (leading zeros, ISO-looking
    # date & time strings that must NOT be reinterpreted as numbers/dates).
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $telp
    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = $tanggal
    $ws.Cells.Item($row, 6).NumberFormat = "@"
    $ws.Cells.Item($row, 6).Value = $jam

    # Jumlah (G) is numeric, except the last row which uses "10 orang".
    if ($jumlah -is [string]) {
        $ws.Cells.Item($row, 7).NumberFormat = "@"
        $ws.Cells.Item($row, 7).Value = $jumlah
    } else {
        $ws.Cells.Item($row, 7).Value = $jumlah
    }

    $ws.Cells.Item($row, 8).Value = $bagian

    # Re-apply the pristine border/general-format style to the whole row
    # without touching the values/types we just wrote.
    $ws.Range("A20:H20").Copy()
    $ws.Range("A" + $row + ":H" + $row).PasteSpecial(-4122)
}

# Row 6: replace the existing "syahrul" guest with the new "Nafisa" entry
Set-GuestRow 6 1 "Nafisa " "PT Pelita Jaya Harapann" "002131618" "2023-06-19" "09:00:28" 23 "Umum"

# New rows 7-12: "Nafisa Azzahra" / "PT Sinar Kasih" visits
Set-GuestRow 7  2 "Nafisa Azzahra" "PT Sinar Kasih" "082233659" "2023-06-19" "09:00:28" 23 "Komisi C"
Set-GuestRow 8  3 "Nafisa Azzahra" "PT Sinar Kasih" "082233659" "2023-06-19" "09:00:27" 20 "Persidangan"
Set-GuestRow 9  4 "Nafisa Azzahra" "PT Sinar Kasih" "082233659" "2023-06-19" "10:00:55" 20 "Persidangan"
Set-GuestRow 10 5 "Nafisa Azzahra" "PT Sinar Kasih" "082233659" "2023-06-17" "00:05:00" 12 "Umum"
Set-GuestRow 11 6 "Nafisa Azzahra" "PT Sinar Kasih" "082233659" "2023-06-16" "22:54:00" 12 "Umum"
Set-GuestRow 12 7 "Nafisa Azzahra" "PT Sinar Kasih" "082233659" "2023-06-01" "03:51:00" "10 orang" "Komisi-A"

# Drop the scratch format-donor row so it doesn't leak into the sheet
$ws.Range("A20:H20").Clear()

# Update the report period text (A4, merged A4:H4)
$ws.Range("A4").Value = "Periode 2022-03-17 - 2023-06-21"

# Selection follows the last row, like Excel does after entering a row
$ws.Range("A12:H12").Select()
